$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("E1").Value = "Completed Courses"
$ws.Range("F1").Value = "Uncompleted Courses"

# Row 2: no completed->uncompleted split needed, just clear the dates column
$ws.Range("F2").Value = "[]"

# Row 3: move two courses from E (completed) to F (uncompleted)
$ws.Range("E3").Value = "['Supervisor Safety Training (2 hrs)', 'Employee Safety (1 hr)']"
$ws.Range("F3").Value = "['Counterintelligence (1 hr)', 'HIPAA Training (1 hr)']"

# Row 4: clear the dates column
$ws.Range("F4").Value = "[]"

# Row 5: clear the dates column
$ws.Range("F5").Value = "[]"
